$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" column (E16:E21) currently lists the six periods in
# descending order (1706,1705,1704,1703,1702,1701). The account-statement
# database was refreshed and the periods are now listed in ascending order
# (1701,1702,1703,1704,1705,1706) so that future statements can keep being
# appended after 1706.
$ws.Range("E16").Value = "1701"
$ws.Range("E17").Value = "1702"
$ws.Range("E18").Value = "1703"
$ws.Range("E19").Value = "1704"
$ws.Range("E20").Value = "1705"
$ws.Range("E21").Value = "1706"
